# Update gh-pages to output generated at 456a3b4
# Refresh the "want to go" counts (column F) and a couple of min-price
# fixes (column G) across the four sheets: 展览, 演出, 本地生活, 全部类型.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型

# 展览
$ws1.Range("F6").Value = 274
$ws1.Range("F7").Value = 13148
$ws1.Range("G7").Value = 85
$ws1.Range("F8").Value = 66
$ws1.Range("F10").Value = 282
$ws1.Range("F11").Value = 4101
$ws1.Range("F12").Value = 6721
$ws1.Range("F13").Value = 62
$ws1.Range("F15").Value = 3516
$ws1.Range("F21").Value = 129
$ws1.Range("F22").Value = 3657
$ws1.Range("F23").Value = 104
$ws1.Range("F25").Value = 3707
$ws1.Range("F26").Value = 3707
$ws1.Range("F27").Value = 422
$ws1.Range("F28").Value = 1920
$ws1.Range("F30").Value = 243
$ws1.Range("F31").Value = 6864
$ws1.Range("F34").Value = 1785
$ws1.Range("F35").Value = 2032
$ws1.Range("F36").Value = 1303
$ws1.Range("F37").Value = 112
$ws1.Range("F38").Value = 1084
$ws1.Range("F44").Value = 1149
$ws1.Range("F46").Value = 145
$ws1.Range("F47").Value = 1230
$ws1.Range("F48").Value = 1836
$ws1.Range("F49").Value = 72
$ws1.Range("F50").Value = 166

# 演出
$ws2.Range("F14").Value = 105

# 本地生活
$ws3.Range("F2").Value = 466
$ws3.Range("F3").Value = 636
$ws3.Range("F4").Value = 29

# 全部类型
$ws4.Range("F6").Value = 466
$ws4.Range("F7").Value = 636
$ws4.Range("F8").Value = 29
$ws4.Range("F9").Value = 274
$ws4.Range("F10").Value = 13148
$ws4.Range("G10").Value = 85
$ws4.Range("F11").Value = 66
$ws4.Range("F14").Value = 282
$ws4.Range("F15").Value = 4101
$ws4.Range("F16").Value = 6721
$ws4.Range("F17").Value = 62
$ws4.Range("F18").Value = 3516
$ws4.Range("F25").Value = 129
$ws4.Range("F28").Value = 3708
$ws4.Range("F29").Value = 422
$ws4.Range("F31").Value = 243
$ws4.Range("F32").Value = 6864
$ws4.Range("F33").Value = 105
$ws4.Range("F36").Value = 1785
$ws4.Range("F37").Value = 2032
$ws4.Range("F38").Value = 1303
$ws4.Range("F39").Value = 112
$ws4.Range("F40").Value = 1084
$ws4.Range("F44").Value = 1149
$ws4.Range("F45").Value = 145
$ws4.Range("F47").Value = 1836
$ws4.Range("F48").Value = 72
$ws4.Range("F50").Value = 166
